$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.192249655723572
$ws.Range("B1").Value = 2.309991359710693
$ws.Range("C1").Value = 6.687001705169678
$ws.Range("D1").Value = 2.32832932472229
$ws.Range("E1").Value = 1.188352465629578
